$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D2').Value = '28.806.62'

# Row 3
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('D3').Value = '1.811.64'

# Row 4
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'

# Row 5
$ws.Range('E5').Value = '  -2.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.44'
$ws.Range('D5').Style = 'Normal'

# Row 6
$ws.Range('E6').Value = '  -3.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5894'
$ws.Range('D6').Style = 'Normal'

# Row 8
$ws.Range('B8').Value = 'WrappedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D8').Value = '1.850.50'

# Row 9
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('E9').Value = '  -2.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2744'
$ws.Range('D9').Style = 'Normal'

# Row 10
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('E10').Value = '  -5.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06720'
$ws.Range('D10').Style = 'Normal'

# Row 11
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('E11').Value = '  -3.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.01'
$ws.Range('D11').Style = 'Normal'

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07492'
$ws.Range('D12').Style = 'Normal'

# Row 13
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.668'
$ws.Range('D13').Style = 'Normal'

# Row 14
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6227'
$ws.Range('D14').Style = 'Normal'

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D15').Value = '2.058.10'

# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E16').Value = '  -10.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009047'
$ws.Range('D16').Style = 'Normal'

# Row 17
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E17').Value = '  -6.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '74.40'
$ws.Range('D17').Style = 'Normal'

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D18').Value = '28.577.71'

# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E19').Value = '  -8.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.437'
$ws.Range('D19').Style = 'Normal'

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('D20').Style = 'Normal'

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E21').Value = '  -10.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '207.36'
$ws.Range('D21').Style = 'Normal'

# Row 22
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('E22').Value = '  -4.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.33'
$ws.Range('D22').Style = 'Normal'

# Row 23
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E23').Value = '  -4.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.757'
$ws.Range('D23').Style = 'Normal'

# Row 24
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D24').Style = 'Normal'

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.14'
$ws.Range('D25').Style = 'Normal'

# Row 26
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E26').Value = '  -3.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.802'
$ws.Range('D26').Style = 'Normal'

# Row 27
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1261'
$ws.Range('D27').Style = 'Normal'

# Row 28
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.28'
$ws.Range('D28').Style = 'Normal'

# Row 29
$ws.Range('E29').Value = '  -4.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.424'
$ws.Range('D29').Style = 'Normal'

# Row 30
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E30').Value = '  -7.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06108'
$ws.Range('D30').Style = 'Normal'

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E31').Value = '  -2.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.423'
$ws.Range('D31').Style = 'Normal'

# Row 32
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.715'
$ws.Range('D32').Style = 'Normal'

# Row 33
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E33').Value = '  -3.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.689'
$ws.Range('D33').Style = 'Normal'

# Row 34
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.698'
$ws.Range('D34').Style = 'Normal'

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E35').Value = '  -7.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.040'
$ws.Range('D35').Style = 'Normal'

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E36').Value = '  -3.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6306'
$ws.Range('D36').Style = 'Normal'

# Row 37
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.519'
$ws.Range('D37').Style = 'Normal'

# Row 38
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E38').Value = '  -1.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.736'
$ws.Range('D38').Style = 'Normal'

# Row 39
$ws.Range('E39').Value = '  -4.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01686'
$ws.Range('D39').Style = 'Normal'

# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E40').Value = '  -4.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.337'
$ws.Range('D40').Style = 'Normal'

# Row 41
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('E41').Value = '  -8.42%  '
$ws.Range('D41').Value = '1.126.89'

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E42').Value = '  -6.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8634'
$ws.Range('D42').Style = 'Normal'

# Row 43
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('D43').Style = 'Normal'

# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.45'
$ws.Range('D44').Style = 'Normal'

# Row 45
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D45').Value = '1.968.91'

# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E46').Value = '  -5.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.08'
$ws.Range('D46').Style = 'Normal'

# Row 47
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E47').Value = '  -3.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000112'
$ws.Range('D47').Style = 'Normal'

# Row 48
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.564'
$ws.Range('D48').Style = 'Normal'

# Row 49
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05466'
$ws.Range('D49').Style = 'Normal'

# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4511'
$ws.Range('D50').Style = 'Normal'

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E51').Value = '  -3.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.254'
$ws.Range('D51').Style = 'Normal'
